$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("re_profiles")

$script:scratchRow = 500

function Move-Range {
    param($ws, [string]$srcAddr, [string]$dstAddr)
    $src = $ws.Range($srcAddr)
    $dst = $ws.Range($dstAddr)

    # Stage through a far-away scratch area so overlapping source/destination
    # ranges never alias each other mid-operation.
    $rowCount = $src.Rows.Count
    $colCount = $src.Columns.Count
    $scratchTopLeft = $ws.Cells.Item($script:scratchRow, 2)
    $scratch = $ws.Range($scratchTopLeft, $ws.Cells.Item($script:scratchRow + $rowCount - 1, 1 + $colCount))
    $script:scratchRow = $script:scratchRow + $rowCount + 2

    $src.Copy($scratch)
    foreach ($cell in $src) {
        $cell.Style = "Normal"
        $cell.Value = $null
    }

    $scratch.Copy($dst)
    foreach ($cell in $scratch) {
        $cell.Style = "Normal"
        $cell.Value = $null
    }
}

# Block B: H9:K11 -> G9:J11 (shift left 1 column)
Move-Range $ws "H9:K11" "G9"

# Block C: O9:R11 -> L9:O11 (shift left 3 columns)
Move-Range $ws "O9:R11" "L9"

# Small lookup table M2:O4 -> Q9:S11
Move-Range $ws "M2:O4" "Q9"

Write-Output "done"
